# The "variables" sheet (sheet index 1) holds the sample table in A1:N19.
# Rows 2-19 are nine samples, each duplicated into two consecutive rows
# (2&3, 4&5, 6&7, ... 18&19). The edit removes the duplicate (second) row
# of every pair, leaving one row per sample and compacting the table to
# A1:N10.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete from the bottom up so earlier row numbers stay valid.
$rowsToDelete = @(19, 17, 15, 13, 11, 9, 7, 5, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
